# "Add files via upload" - Product Backlog update
# Rewrites the "PB Items" sheet (user stories/features table) with the
# refreshed backlog content, adjusts a few column widths, and moves the
# active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PB Items")

# --- Row 6 : US1 / F1 ---------------------------------------------------
$ws.Range("C6").Value = "customer"
$ws.Range("D6").Value = "select different comics from a list (search by categories, title)"
$ws.Range("E6").Value = "..I can buy them"
$ws.Range("I6").Value = "Jonas"
$ws.Range("K6").Value = "User sees a list of different comics and can select one to see it in a single view, there is also a search bar"

# --- Row 7 : US2 / F2 ---------------------------------------------------
$ws.Range("C7").Value = "customer"
$ws.Range("D7").Value = "create an account"
$ws.Range("E7").Value = "..I can pay online"
$ws.Range("I7").Value = "Jonas"
$ws.Range("K7").Value = "User can deposit his name, an address, an Email-Address, and a Paypal-Account"

# --- Row 8 : US3 / F3 ---------------------------------------------------
$ws.Range("C8").Value = "salesman"
$ws.Range("D8").Value = "create an account with email confirmation"
$ws.Range("E8").Value = "..I get messages when somebody buys something"
$ws.Range("I8").Value = "Jonas"
$ws.Range("K8").Value = "User can deposit the same informations as in other accounts, but also see the informations of customer-accounts which are hidden for normal users"

# --- Row 9 : US4 / F4 ---------------------------------------------------
$ws.Range("C9").Value = "customer, salesman"
$ws.Range("D9").Value = "login"
$ws.Range("E9").Value = "..I can access my account"
$ws.Range("F9").Value = "Must have"
$ws.Range("I9").Value = "Jonas"
$ws.Range("K9").Value = "User can access the shop and access an account-view"

# --- Row 10 : US5 / F5 ---------------------------------------------------
$ws.Range("C10").Value = "customer "
$ws.Range("D10").Value = "order a comic book"
$ws.Range("E10").Value = "..it can be sent to my address"
$ws.Range("F10").Value = "Must Have"
$ws.Range("I10").Value = "Jonas"
$ws.Range("K10").Value = "User can press a buy-button in the single view of a comic-object, then it will be sent to his address"

# --- Row 11 : US6 / F6 ---------------------------------------------------
$ws.Range("C11").Value = "salesman"
$ws.Range("D11").Value = "add new comics to the list"
$ws.Range("E11").Value = "..I can sell them"
$ws.Range("F11").Value = "Must Have"
$ws.Range("I11").Value = "Jonas"
$ws.Range("K11").Value = "User can create a new comic-object"

# --- Row 12 : US7 / F7 (new row) ----------------------------------------
$ws.Range("A12").Value = "US7"
$ws.Range("B12").Value = "F7"
$ws.Range("C12").Value = "external user"
$ws.Range("D12").Value = "have a restful API"
$ws.Range("E12").Value = "..I can retrieve data about a comic book"
$ws.Range("F12").Value = "Must Have"
$ws.Range("G12").Value = 7
$ws.Range("I12").Value = "Jonas"
$ws.Range("K12").Value = "User can make CRUD-Operations on a comic-object through an API"

# --- Row 13 : US8 / F8 (new row) ----------------------------------------
$ws.Range("A13").Value = "US8"
$ws.Range("B13").Value = "F8"
$ws.Range("C13").Value = "salesman"
$ws.Range("D13").Value = "create a newsletter"
$ws.Range("E13").Value = "..I can give it to users"
$ws.Range("F13").Value = "Could Have"
$ws.Range("G13").Value = 8
$ws.Range("I13").Value = "Jonas"
$ws.Range("K13").Value = "User can save a list of comic-objects and announcement-objects and transfer it to a Email-Address"

# --- Row 14 : US9 / F9 (new row) ----------------------------------------
$ws.Range("A14").Value = "US9"
$ws.Range("B14").Value = "F9"
$ws.Range("C14").Value = "writer"
$ws.Range("D14").Value = "add announcements"
$ws.Range("E14").Value = "..they are advertisement"
$ws.Range("F14").Value = "Could Have"
$ws.Range("G14").Value = 9
$ws.Range("I14").Value = "Jonas"
$ws.Range("K14").Value = "User can create special announcement-objects"

# --- Rows 15/16 : old US10 / US11 rows are now blank ---------------------
$ws.Range("A15").ClearContents()
$ws.Range("A16").ClearContents()

# --- Column widths --------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 17.333333333333336
$ws.Columns.Item(4).ColumnWidth = 50.83333333333333
$ws.Columns.Item(5).ColumnWidth = 43.5
$ws.Columns.Item(9).ColumnWidth = 13.0

# --- Selection -------------------------------------------------------------
$ws.Range("K11").Select() | Out-Null
